$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1273.6533203125
$ws.Range("C2").Value = 0.9555
$ws.Range("D2").Value = 0.9143000245094299
$ws.Range("E2").Value = 1.504899978637695
$ws.Range("F2").Value = 0.713699996471405
$ws.Range("H2").Value = 0.9263

$ws.Range("B3").Value = 1190.009155273438
$ws.Range("C3").Value = 0.952
$ws.Range("D3").Value = 0.9243
$ws.Range("E3").Value = 1.851799964904785
$ws.Range("F3").Value = 0.7566999793052673
$ws.Range("H3").Value = 1.0156

$ws.Range("B4").Value = 712.0230102539062
$ws.Range("C4").Value = 0.8318
$ws.Range("D4").Value = 0.8242
$ws.Range("E4").Value = 1.070299983024597
$ws.Range("F4").Value = 0.7006000280380249
$ws.Range("H4").Value = 0.1284

$ws.Range("B5").Value = 724.401123046875
$ws.Range("C5").Value = 0.7641
$ws.Range("D5").Value = 0.7577
$ws.Range("E5").Value = 1.19159996509552
$ws.Range("F5").Value = 0.6705999970436096
$ws.Range("H5").Value = -0.4607

$ws.Range("B6").Value = 1158.219482421875
$ws.Range("C6").Value = 0.9163
$ws.Range("D6").Value = 0.9139
$ws.Range("E6").Value = 1.238100051879883
$ws.Range("F6").Value = 0.7195000052452087
$ws.Range("H6").Value = 0.9231

$ws.Range("B7").Value = 897.79638671875
$ws.Range("C7").Value = 0.9041
$ws.Range("D7").Value = 0.8964999914169312
$ws.Range("E7").Value = 1.225800037384033
$ws.Range("F7").Value = 0.7687000036239624
$ws.Range("H7").Value = 0.7692

$ws.Range("B8").Value = 1001.961303710938
$ws.Range("C8").Value = 0.8978
$ws.Range("D8").Value = 0.8872
$ws.Range("E8").Value = 1.255599975585938
$ws.Range("F8").Value = 0.7868000268936157
$ws.Range("H8").Value = 0.6862

$ws.Range("B9").Value = 6958.06396484375
$ws.Range("C9").Value = 0.8966
$ws.Range("D9").Value = 0.8848
$ws.Range("E9").Value = 1.851799964904785
$ws.Range("F9").Value = 0.6705999970436096
$ws.Range("H9").Value = 3.9881
